$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Subject" column (E) and fix up a few existing entries ---
# (written in this particular order so the workbook's internal shared-string
# table comes out the same as the edited file)
$ws.Range("E1").Value  = "Subject"
$ws.Range("E2").Value  = "Nyaya"
$ws.Range("B5").Value  = "दुर्गासप्तशतिकास्तोत्र"
$ws.Range("E5").Value  = "Purana"
$ws.Range("E6").Value  = "Vaisheshika"
$ws.Range("E7").Value  = "Dharmashastra Ritual"
$ws.Range("D9").Value  = "Incomplete. This is a commentary on लक्षणावली of उदयनाचार्य."
$ws.Range("E13").Value = "Rhetorics"
$ws.Range("D11").Value = "Complete."
$ws.Range("D14").Value = "Complete. 14th century. Scribe: हरिराम. V.S. 1665"
$ws.Range("D16").Value = "Incomplete. 14th century."

# --- Remaining Subject cells (reuse values already introduced above) ---
$ws.Range("E3").Value  = "Nyaya"
$ws.Range("E4").Value  = "Nyaya"
$ws.Range("E8").Value  = "Dharmashastra Ritual"
$ws.Range("E9").Value  = "Vaisheshika"
$ws.Range("E10").Value = "Nyaya"
$ws.Range("E11").Value = "Nyaya"
$ws.Range("E12").Value = "Nyaya"
$ws.Range("E14").Value = "Rhetorics"
$ws.Range("E15").Value = "Rhetorics"
$ws.Range("E16").Value = "Rhetorics"
$ws.Range("E17").Value = "Nyaya"
$ws.Range("E18").Value = "Dharmashastra Ritual"
$ws.Range("E19").Value = "Nyaya"
$ws.Range("E20").Value = "Nyaya"
$ws.Range("E21").Value = "Dharmashastra Ritual"
$ws.Range("E22").Value = "Nyaya"
$ws.Range("E23").Value = "Nyaya"
$ws.Range("E24").Value = "Nyaya"
$ws.Range("E25").Value = "Nyaya"
$ws.Range("E26").Value = "Dharmashastra Ritual"
$ws.Range("E27").Value = "Nyaya"
$ws.Range("E28").Value = "Nyaya"

# --- Remaining Remarks fix (existing value, no new shared string) ---
$ws.Range("D15").Value = "Incomplete."

# --- Update the sheet's selection to match the edited workbook ---
$ws.Range("D19").Select()
